$d = $word.ActiveDocument

$replacements = @(
    @("109×4=", "493×3="),
    @("553×3=", "209×3="),
    @("453×9=", "994×3="),
    @("123×4=", "945×7="),
    @("730×3=", "809×2="),
    @("507×3=", "690×9="),
    @("394×2=", "757×3="),
    @("895×5=", "533×5="),
    @("899×6=", "443×2="),
    @("680×3=", "684×6="),
    @("633×4=", "829×4="),
    @("599×6=", "129×2="),
    @("784×6=", "817×5="),
    @("347×3=", "926×9="),
    @("232×6=", "607×4="),
    @("229×4=", "614×4="),
    @("337×2=", "478×5="),
    @("592×9=", "475×2="),
    @("816×5=", "519×7="),
    @("721×9=", "665×4="),
    @("533×4=", "837×9="),
    @("563×2=", "842×5="),
    @("498×8=", "960×4="),
    @("148×8=", "278×2="),
    @("130×2=", "952×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
